# TC03_Canine_StudyUBC02-Breed_Diagnosis_PrimDiseaseSite.xlsx
# "updated ubc2 10 scripts, renamed test suites with w, commiting stashed changes"
#
# The only functional change is in the "startup" sheet: the Cypher query
# stored in cell B2 (the "CasesTab" query) dropped its trailing
# coalesce(co.cohort_description, ...) AS Cohort projection - the
# Cohort column is no longer pulled back for the Cases tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")
$query = $cell.Value2

$needle = ",`r`n        coalesce(co.cohort_description, '') AS ``Cohort``"
$needle2 = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"

if ($query.Contains($needle)) {
    $query = $query.Replace($needle, "")
} elseif ($query.Contains($needle2)) {
    $query = $query.Replace($needle2, "")
}

$cell.Value2 = $query

# Re-saving this workbook (newer Excel build) re-measured the wrapped,
# autofit row heights for the three wrapped query rows - the text box
# widths didn't move, but the recalculated font metrics trimmed each
# row a little.
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 261
